$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 (pushes "Description" and everything below it down by one)
$ws.Rows.Item(11).Insert()

# Match the look of the surrounding data rows (style/format) on the newly inserted row
$ws.Cells.Item(10, 1).Resize(1, 2).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new "Jurisdiction" property row (value left blank, as in the source IG)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the build Date value (row 8 - unaffected by the row-11 insert)
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"
